# edit.ps1 - applies the documented changes to 323 Assignment 2 Documentation.docx
#
# Strategy: for each paragraph that needs its run content restructured
# (including inserting/removing <w:proofErr/> markers and bookmarks), we
# capture the paragraph's Range (which spans from the start of the
# paragraph up to -- but not including -- the paragraph mark) and call
# Range.InsertXML() with a full replacement <w:p> (including the original
# <w:pPr/>) wrapped in the package/part XML that Word expects. This
# reliably replaces the paragraph's contents in a single atomic step and
# avoids quirks in the hosted object model around partial-run editing.

$d = $word.ActiveDocument

function Set-ParagraphXml {
    param(
        [int]$Index,
        [string]$InnerXml
    )
    $para = $d.Paragraphs($Index)
    $range = $para.Range
    $pkg = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>" + $InnerXml + "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    $range.InsertXML($pkg)
}

# ---------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that sits right before
#    "The second assignment is to write a syntax analyzer...".
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2) "...an text file. So far "test1.txt" is only test with the value
#    "b+c+d"" -> "...an text file. "test1.txt" and "test2.txt" were
#    tested in the process"
#    (paragraph beginning "This program is programmed on VisualStudio 2017...")
# ---------------------------------------------------------------------
Set-ParagraphXml 10 @'
<w:p><w:pPr><w:pStyle w:val="Heading1"/><w:tabs><w:tab w:val="left" w:pos="820"/></w:tabs><w:spacing w:line="480" w:lineRule="auto"/><w:ind w:left="0" w:firstLine="0"/><w:rPr><w:b w:val="0"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="0"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">This program is programmed on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b w:val="0"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>VisualStudio</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b w:val="0"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> 2017. When running the program, it will ask for </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b w:val="0"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>an</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b w:val="0"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> text file.</w:t></w:r><w:r><w:rPr><w:b w:val="0"/><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> &#8220;test1.txt&#8221; and &#8220;test2.txt&#8221; were tested in the process</w:t></w:r></w:p>
'@

# ---------------------------------------------------------------------
# 3) "The lexemes are then runned into..." -> "The lexemes goes into..."
#    plus two new sentences about the syntax analyzer's starting symbol
#    and predictive tables.
#    (paragraph beginning "The design of our program uses our lexical
#    analyzer program...")
# ---------------------------------------------------------------------
Set-ParagraphXml 12 @'
<w:p><w:pPr><w:pStyle w:val="BodyText"/><w:spacing w:before="8" w:line="480" w:lineRule="auto"/><w:ind w:right="386"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">The design of our program uses our lexical analyzer program to parse the source file into tokens and lexemes. The lexemes </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>goes</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> into the syntax analyzer to check if the format is correct.</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> The syntax analyzer implements its own starting </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>symbol</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> so we will not be using the &#8220;$$&#8221; or &#8220;%%&#8221;. The syntax analyzer uses a </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>predictive tables</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> using the first and follow set.</w:t></w:r></w:p>
'@

# ---------------------------------------------------------------------
# 4) "So far the program can only handle a correct syntax with limited
#    characters. There is no handling for errors yet." ->
#    "When reading the text file during the lexical analysis, it does
#    not detect the last character of the file. So to compensate for
#    this error is to have an extra space or a new line at the end of
#    the text file."
#    (the "Any Limitation" body paragraph)
# ---------------------------------------------------------------------
Set-ParagraphXml 14 @'
<w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="1635"/></w:tabs><w:spacing w:before="3" w:line="480" w:lineRule="auto"/><w:ind w:right="386"/><w:rPr><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">When reading the text file during the lexical analysis, it does not detect the last character of the file. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>So</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:i/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> to compensate for this error is to have an extra space or a new line at the end of the text file.</w:t></w:r></w:p>
'@

# ---------------------------------------------------------------------
# 5) Delete the whole "So far this is a first iteration, will update
#    when final iteration come." paragraph.
# ---------------------------------------------------------------------
$d.Paragraphs(16).Range.Delete()

# ---------------------------------------------------------------------
# 6) Populate the (now) empty BodyText paragraph that follows with the
#    new closing sentence, and move the "_GoBack" bookmark to sit right
#    after it.
# ---------------------------------------------------------------------
Set-ParagraphXml 16 @'
<w:p><w:pPr><w:pStyle w:val="BodyText"/><w:spacing w:line="480" w:lineRule="auto"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>This syntax analysis did not implement the optional implementations such as declaration, if-else, while, function</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
